# Weekly update: insert one new record row for "Choclo" (Vega Central
# Mapocho de Santiago) before the existing row 187, shifting every
# following row down by one (A1:R291 -> A1:R292).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 187; Excel shifts rows 187:291 down to 188:292
# and grows the sheet dimension automatically.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new weekly record.
$ws.Cells.Item(187, 1).Value  = 9
$ws.Cells.Item(187, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(187, 3).Value  = "Metropolitana"
$ws.Cells.Item(187, 4).Value  = 44460
$ws.Cells.Item(187, 5).Value  = 13
$ws.Cells.Item(187, 6).Value  = 100112024
$ws.Cells.Item(187, 7).Value  = "Choclo"
$ws.Cells.Item(187, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(187, 9).Value  = "Primera"
$ws.Cells.Item(187, 10).Value = 25
$ws.Cells.Item(187, 11).Value = 33000
$ws.Cells.Item(187, 12).Value = 34000
$ws.Cells.Item(187, 13).Value = 33520
$ws.Cells.Item(187, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(187, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(187, 16).Value = 479
$ws.Cells.Item(187, 17).Value = 70
$ws.Cells.Item(187, 18).Value = "Hortaliza"
